$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2160193333333333
$ws.Range("H2").Value = 0.648058
$ws.Range("I2").Value = 0.02486881244588016
$ws.Range("J2").Value = 0.02486881244588016
$ws.Range("M2").Value = 0.07919566666666666
$ws.Range("N2").Value = 0.237587
$ws.Range("O2").Value = 0.08232403487459106
$ws.Range("P2").Value = 0.08232403487459104
$ws.Range("Q2").Value = 0.01710779511622222
$ws.Range("R2").Value = 0.153970156046
$ws.Range("S2").Value = 0.002047300983084302
$ws.Range("T2").Value = 0.002047300983084302
$ws.Range("G3").Value = 0.2160193333333333
$ws.Range("H3").Value = 0.648058
$ws.Range("I3").Value = 0.02486881244588016
$ws.Range("J3").Value = 0.02486881244588016
$ws.Range("O3").Value = 0.9044794902837771
$ws.Range("P3").Value = 0.9044794902837769
$ws.Range("Q3").Value = 0.1879602941008889
$ws.Range("R3").Value = 1.691642646908
$ws.Range("S3").Value = 0.02249333080501253
$ws.Range("T3").Value = 0.02249333080501254
$ws.Range("G4").Value = 0.2160193333333333
$ws.Range("H4").Value = 0.648058
$ws.Range("I4").Value = 0.02486881244588016
$ws.Range("J4").Value = 0.02486881244588016
$ws.Range("M4").Value = 0.012695
$ws.Range("N4").Value = 0.038085
$ws.Range("O4").Value = 0.01319647484163191
$ws.Range("P4").Value = 0.01319647484163191
$ws.Range("Q4").Value = 0.002742365436666667
$ws.Range("R4").Value = 0.02468128893
$ws.Range("S4").Value = 0.00032818065778332
$ws.Range("T4").Value = 0.0003281806577833201
$ws.Range("I5").Value = 0.9551554900377276
$ws.Range("J5").Value = 0.9551554900377278
$ws.Range("M5").Value = 0.07919566666666666
$ws.Range("N5").Value = 0.237587
$ws.Range("O5").Value = 0.08232403487459106
$ws.Range("P5").Value = 0.08232403487459104
$ws.Range("Q5").Value = 0.6570721647147777
$ws.Range("R5").Value = 5.913649482432999
$ws.Range("S5").Value = 0.078632253872523
$ws.Range("T5").Value = 0.078632253872523
$ws.Range("I6").Value = 0.9551554900377276
$ws.Range("J6").Value = 0.9551554900377278
$ws.Range("O6").Value = 0.9044794902837771
$ws.Range("P6").Value = 0.9044794902837769
$ws.Range("S6").Value = 0.8639185507710752
$ws.Range("T6").Value = 0.8639185507710753
$ws.Range("I7").Value = 0.9551554900377276
$ws.Range("J7").Value = 0.9551554900377278
$ws.Range("M7").Value = 0.012695
$ws.Range("N7").Value = 0.038085
$ws.Range("O7").Value = 0.01319647484163191
$ws.Range("P7").Value = 0.01319647484163191
$ws.Range("Q7").Value = 0.1053281256683333
$ws.Range("R7").Value = 0.9479531310149999
$ws.Range("S7").Value = 0.01260468539412947
$ws.Range("T7").Value = 0.01260468539412947
$ws.Range("G8").Value = 0.173491
$ws.Range("H8").Value = 0.520473
$ws.Range("I8").Value = 0.01997281943922393
$ws.Range("J8").Value = 0.01997281943922393
$ws.Range("M8").Value = 0.07919566666666666
$ws.Range("N8").Value = 0.237587
$ws.Range("O8").Value = 0.08232403487459106
$ws.Range("P8").Value = 0.08232403487459104
$ws.Range("Q8").Value = 0.01373973540566666
$ws.Range("R8").Value = 0.123657618651
$ws.Range("S8").Value = 0.001644243084058581
$ws.Range("T8").Value = 0.001644243084058581
$ws.Range("G9").Value = 0.173491
$ws.Range("H9").Value = 0.520473
$ws.Range("I9").Value = 0.01997281943922393
$ws.Range("J9").Value = 0.01997281943922393
$ws.Range("O9").Value = 0.9044794902837771
$ws.Range("P9").Value = 0.9044794902837769
$ws.Range("Q9").Value = 0.1509560226886666
$ws.Range("R9").Value = 1.358604204198
$ws.Range("S9").Value = 0.01806500554591917
$ws.Range("T9").Value = 0.01806500554591918
$ws.Range("G10").Value = 0.173491
$ws.Range("H10").Value = 0.520473
$ws.Range("I10").Value = 0.01997281943922393
$ws.Range("J10").Value = 0.01997281943922393
$ws.Range("M10").Value = 0.012695
$ws.Range("N10").Value = 0.038085
$ws.Range("O10").Value = 0.01319647484163191
$ws.Range("P10").Value = 0.01319647484163191
$ws.Range("Q10").Value = 0.002202468245
$ws.Range("R10").Value = 0.019822214205
$ws.Range("S10").Value = 0.0002635708092461753
$ws.Range("T10").Value = 0.0002635708092461754
$ws.Range("G11").Value = [double]"2.5E-05"
$ws.Range("H11").Value = [double]"7.499999999999999E-05"
$ws.Range("I11").Value = [double]"2.878077168156263E-06"
$ws.Range("J11").Value = [double]"2.878077168156264E-06"
$ws.Range("M11").Value = 0.07919566666666666
$ws.Range("N11").Value = 0.237587
$ws.Range("O11").Value = 0.08232403487459106
$ws.Range("P11").Value = 0.08232403487459104
$ws.Range("Q11").Value = [double]"1.979891666666667E-06"
$ws.Range("R11").Value = [double]"1.7819025E-05"
$ws.Range("S11").Value = [double]"2.369349251630604E-07"
$ws.Range("T11").Value = [double]"2.369349251630605E-07"
$ws.Range("G12").Value = [double]"2.5E-05"
$ws.Range("H12").Value = [double]"7.499999999999999E-05"
$ws.Range("I12").Value = [double]"2.878077168156263E-06"
$ws.Range("J12").Value = [double]"2.878077168156264E-06"
$ws.Range("O12").Value = 0.9044794902837771
$ws.Range("P12").Value = 0.9044794902837769
$ws.Range("Q12").Value = [double]"2.175271666666666E-05"
$ws.Range("R12").Value = 0.0001957744499999999
$ws.Range("S12").Value = [double]"2.603161770051353E-06"
$ws.Range("T12").Value = [double]"2.603161770051354E-06"
$ws.Range("G13").Value = [double]"2.5E-05"
$ws.Range("H13").Value = [double]"7.499999999999999E-05"
$ws.Range("I13").Value = [double]"2.878077168156263E-06"
$ws.Range("J13").Value = [double]"2.878077168156264E-06"
$ws.Range("M13").Value = 0.012695
$ws.Range("N13").Value = 0.038085
$ws.Range("O13").Value = 0.01319647484163191
$ws.Range("P13").Value = 0.01319647484163191
$ws.Range("Q13").Value = [double]"3.17375E-07"
$ws.Range("R13").Value = [double]"2.856375E-06"
$ws.Range("S13").Value = [double]"3.798047294184933E-08"
$ws.Range("T13").Value = [double]"3.798047294184934E-08"
